$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 104
$ws.Range("A104").Value = 'Psychiatric comorbidity and social adjustment difficulties in children with disruptive mood dysregulation disorder: A national epidemiological study'
$ws.Range("B104").Value = 'Yu-Ju Lin, Wan-Ling Tseng, Susan Shur-Fen Gau'
$ws.Range("C104").Value = 'Journal of Affective Disorders'
$ws.Range("D104").Value = 2021
$ws.Range("E104").Value = 'https://doi.org/10.1016/j.jad.2020.12.039'
$ws.Range("F104").Value = 'Background
We examined correlates, psychiatric comorbidities, and social adjustment difficulties in children with disruptive mood dysregulation disorder (DMDD) from a national epidemiological study in Taiwan.
Methods
The sample consisted of 4816 children, who were 3rd, 5th, and 7th graders from schools randomly chosen based on the urbanization level in a recent national survey of childhood mental disorders. Among the 4816 children (2520 boys, 52.3%) interviewed using the Kiddie epidemiologic version of the Schedule for Affective Disorders and Schizophrenia (K-SADS-E) for the DSM-5, 30 children were diagnosed as DMDD (23 boys, 76.7%). They and their parents also reported on the Social Adjustment Inventory for Children and Adolescents (SAICA). We conducted regressions for survey data that controlled for stratification and clustering.
Results
The weighted prevalence of DMDD was 0.3~0.76% in Taiwanese children. Lower parental educational levels, male predominance, higher psychiatric comorbidities, and worse self-report school functions (e.g., more behavioral problems with peers) were observed in children with DMDD than those without. Additional analyses revealed that oppositional defiant disorder (ODD) but not DMDD was related to conduct disorder. Children with ODD with or without DMDD had more problems regarding attitudes toward school, academic performance, and parent-child interaction at home than those with DMDD-only.
Limitations
Small sample size of DMDD.
Conclusions
DMDD is a rare disorder in the community. Children with DMDD had more psychiatric comorbidities, and subjectively experienced more difficulties than those without. DMDD and ODD both resulted in severe impairment yet in different domains.'
$ws.Range("G104").Value = 'ARI referenced in paper'
$ws.Range("H104").Value = 'Observational'

# Row 105
$ws.Range("A105").Value = 'Phasic Versus Tonic Irritability: Differential Associations With Attention-Deficit/Hyperactivity Disorder Symptoms'
$ws.Range("B105").Value = 'Elise M. Cardinale, Gabrielle F. Freitag, Melissa A. Brotman, Daniel S. Pine, Ellen Leibenluft, Katharina Kircanski'
$ws.Range("C105").Value = 'Journal of the American Academy of Child & Adolescent Psychiatry'
$ws.Range("D105").Value = 2021
$ws.Range("E105").Value = 'https://doi.org/10.1016/j.jaac.2020.11.022'
$ws.Range("F105").Value = 'Objective
Irritability is a multifaceted construct in pediatric psychopathology. It has been conceptualized as having a ‘phasic’ dimension and a ‘tonic’ dimension. Disruptive mood dysregulation disorder (DMDD) is defined by the presence of both dimensions. Severe irritability, or DMDD, is highly comorbid with attention-deficit/hyperactivity disorder (ADHD). However, it is unknown whether the presence of ADHD modulates the expression of phasic and tonic irritability.
Method
A data-driven, latent variable approach was used to examine irritability and ADHD symptoms in a transdiagnostic pediatric sample (N=489) with primary DMDD, ADHD, subclinical irritability symptoms, or no diagnosis. Using latent profile analyses (LPA), we identified four classes: high levels of both irritability and ADHD symptoms, high levels of irritability and moderate levels of ADHD symptoms, moderate levels of irritability and high levels of ADHD symptoms, and low levels of both irritability and ADHD symptoms. Confirmatory factor analysis operationalized phasic irritability and tonic irritability.
Results
As expected, the two latent classes characterized by high overall irritability exhibited the highest levels of both phasic and tonic irritability. However, between these two high irritability classes, highly comorbid ADHD symptoms were associated with significantly greater phasic irritability than were moderately comorbid ADHD symptoms. In contrast, the two high irritability groups did not differ on levels of tonic irritability.
Conclusion
These findings suggest that phasic, but not tonic, irritability has a significant association with ADHD symptoms, and that phasic and tonic might be distinct, though highly related, irritability dimensions. Future research should investigate potential mechanisms underlying this differential association.'
$ws.Range("G105").Value = 'ARI used as part of measures'
$ws.Range("H105").Value = 'Observational'

# Row 106
$ws.Range("A106").Value = 'Neural mechanisms of reward processing in adolescent irritability'
$ws.Range("B106").Value = 'Maria Kryza‐Lacombe, Brianna Hernandez, Cassidy Owen, Richard C. Reynolds, Lauren S. Wakschlag, Lea R. Dougherty, Jillian L. Wiggins'
$ws.Range("C106").Value = 'Developmental Psychobiology'
$ws.Range("D106").Value = 2021
$ws.Range("E106").Value = 'https://doi.org/10.1002/dev.22090'
$ws.Range("F106").Value = 'Irritability is impairing and prevalent across pediatric psychiatric disorders and typical development, yet its neural mechanisms are largely unknown. This study evaluated the relation between adolescent irritability and reward‐related brain function as a candidate neural mechanism. Adolescents from intervention‐seeking families in the community (N = 52; mean age = 13.80, SD = 1.94) completed a monetary incentive delay task to assess reward anticipation and feedback (reward receipt and omission) during fMRI acquisition. Whole‐brain analyses, controlling for age, examined brain activation and striatal and amygdala connectivity in relation to irritability. Irritability was measured using the parent‐ and youth‐reported Affective Reactivity Index. Irritability was associated with altered reward processing‐related activation and connectivity in multiple networks during reward anticipation and feedback, including increased striatal activation and altered ventral striatum connectivity with prefrontal areas. Our findings suggest that irritability is associated with altered neural patterns during reward processing and that aberrant prefrontal cortex‐mediated top‐down control may be related to irritability. These findings inform our understanding of the etiology of youth irritability and the development of mechanism‐based interventions.'
$ws.Range("G106").Value = 'ARI used as part of measures'
$ws.Range("H106").Value = 'Observational'

# Row 107
$ws.Range("A107").Value = 'Functional connectivity during frustration: a preliminary study of predictive modeling of irritability in youth'
$ws.Range("B107").Value = 'Dustin Scheinost, Javid Dadashkarimi, Emily S. Finn, Caroline G. Wambach, Caroline MacGillivray, Alexandra L. Roule, Tara A. Niendam, Daniel S. Pine, Melissa A. Brotman, Ellen Leibenluft, Wan-Ling Tseng'
$ws.Range("C107").Value = 'Neuropsychopharmacology'
$ws.Range("D107").Value = 2021
$ws.Range("E107").Value = 'https://doi.org/10.1038/s41386-020-00954-8'
$ws.Range("F107").Value = 'Irritability cuts across many pediatric disorders and is a common presenting complaint in child psychiatry; however, its neural mechanisms remain unclear. One core pathophysiological deficit of irritability is aberrant responses to frustrative nonreward. Here, we conducted a preliminary fMRI study to examine the ability of functional connectivity during frustrative nonreward to predict irritability in a transdiagnostic sample. This study included 69 youths (mean age = 14.55 years) with varying levels of irritability across diagnostic groups: disruptive mood dysregulation disorder (n = 20), attention-deficit/hyperactivity disorder (n = 14), anxiety disorder (n = 12), and controls (n = 23). During fMRI, participants completed a frustrating cognitive flexibility task. Frustration was evoked by manipulating task difficulty such that, on trials requiring cognitive flexibility, “frustration” blocks had a 50% error rate and some rigged feedback, while “nonfrustration” blocks had a 10% error rate. Frustration and nonfrustration blocks were randomly interspersed. Child and parent reports of the affective reactivity index were used as dimensional measures of irritability. Connectome-based predictive modeling, a machine learning approach, with tenfold cross-validation was conducted to identify networks predicting irritability. Connectivity during frustration (but not nonfrustration) blocks predicted child-reported irritability (ρ = 0.24, root mean square error = 2.02, p = 0.03, permutation testing, 1000 iterations, one-tailed). Results were adjusted for age, sex, medications, motion, ADHD, and anxiety symptoms. The predictive networks of irritability were primarily within motor-sensory networks; among motor-sensory, subcortical, and salience networks; and between these networks and frontoparietal and medial frontal networks. This study provides preliminary evidence that individual differences in irritability may be associated with functional connectivity during frustration, a phenotype-relevant state.'
$ws.Range("G107").Value = 'ARI used as part of measures'
$ws.Range("H107").Value = 'Observational'
